$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: "Check main menu Home category is functional" ---
# Fix duplicate step numbering (3. -> 4.) and flesh out the expected result.
$ws.Range("E7").Value = @"
1. Open https://abantecart.codifyme.co.nz
2. Check if HOME Menu link is working
3. Mouseover HOME Menu
4. Check if all links in the HOME Menu dropdown are working
"@

$ws.Range("F7").Value = @"
2. HOME Menu button link should be working.
should be working.
3. Dropdown appears .
4. All dropdown links should be working.
"@

# --- Row 8: "Check main menu Apparel & Accessories category is functional" ---
# This test script was previously blank/placeholder-like; give it the same
# look (font colour / wrap formatting) as the other freshly authored rows
# above it (2-7) by copying their format down, then fill in the real steps.
$ws.Range("A7:F7").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows(8).RowHeight = 90

$ws.Range("E8").Value = @"
1. Open https://abantecart.codifyme.co.nz
2. Check if APPAREL & ACCESSORIES Menu link is working
3. Mouseover APPAREL & ACCESSORIES Menu
4. Check if all links in the APPAREL & ACCESSORIES Menu dropdown are working
"@

$ws.Range("F8").Value = @"
2. APPAREL & ACCESSORIES Menu button link should be working.
3. Dropdown appears .
4. All dropdown links should be working.
"@

# --- Row 9: "Check main menu Makeup category is functional" ---
# The leftover placeholder test steps text doesn't belong here; clear it
# and let the row height shrink back to the sheet default.
$ws.Range("E9").ClearContents()
$ws.Rows(9).AutoFit()

# Move the active selection to reflect where the edit was made.
$null = $ws.Range("E18").Select()
